$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the old "Kompetensi Keahlian / Tempat Prakerin" mini-table
#    that used to live at K12:L17 - it is being relocated next to the
#    main table (columns L:M, rows 2-7).
# ------------------------------------------------------------------
$ws.Rows("12:17").Delete()

# ------------------------------------------------------------------
# 2) Add the relocated columns L (Kompetensi Keahlian) and
#    M (Tempat Prakerin) right after the existing table (col K).
#    Copy formatting from analogous existing cells first, then set
#    the values.
# ------------------------------------------------------------------

# Header row (row 2) - same bold/bordered/centered style as K2.
$ws.Range("K2").Copy()
$ws.Range("L2:M2").PasteSpecial(-4122)
$ws.Range("L2").Value = "Kompetensi Keahlian"
$ws.Range("M2").Value = "Tempat Prakerin"

# Data rows (3-7) - same bordered/centered style as the rest of the row.
$ws.Range("A3").Copy()
$ws.Range("L3:M7").PasteSpecial(-4122)

$ws.Range("L3").Value = "Rekayasa Perangkat Lunak"
$ws.Range("L4").Value = "Rekayasa Perangkat Lunak"
$ws.Range("L5").Value = "Rekayasa Perangkat Lunak"
$ws.Range("L6").Value = "Rekayasa Perangkat Lunak"
$ws.Range("L7").Value = "Rekayasa Perangkat Lunak"

$ws.Range("M3").Value = "BKPSDM"
$ws.Range("M4").Value = "CAPIL"
$ws.Range("M5").Value = "KESBANGPOL"
$ws.Range("M6").Value = "BAPPEDA"
$ws.Range("M7").Value = "PBN"

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 3) Column widths: split the old J:K (10-11) width block and give
#    the new L/M (11/13) columns their own widths.
# ------------------------------------------------------------------
$ws.Columns(11).ColumnWidth = 16.3
$ws.Columns(13).ColumnWidth = 24.6

# ------------------------------------------------------------------
# 4) Update the active cell / selection to match where the author
#    left off editing.
# ------------------------------------------------------------------
$ws.Range("J13").Select() | Out-Null
